$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0049
$ws.Range("E2").Value = 0.26
$ws.Range("D3").Value = 0.00734
$ws.Range("E3").Value = 0.26
$ws.Range("D4").Value = 0.0406
$ws.Range("E4").Value = 0.848
$ws.Range("D5").Value = 0.08260000000000001
$ws.Range("E5").Value = 0.848
$ws.Range("D6").Value = 0.0863
$ws.Range("E6").Value = 0.848
$ws.Range("A7").Formula = '=HYPERLINK("pathways/Lysine_Degradation.csv","Lysine Degradation")'
$ws.Range("B7").Value = 5
$ws.Range("D7").Value = 0.08989999999999999
$ws.Range("E7").Value = 0.848
$ws.Range("A8").Formula = '=HYPERLINK("pathways/Alanine_Metabolism.csv","Alanine Metabolism")'
$ws.Range("B8").Value = 9
$ws.Range("D8").Value = 0.0901
$ws.Range("E8").Value = 0.848
$ws.Range("D9").Value = 0.103
$ws.Range("E9").Value = 0.848
$ws.Range("D10").Value = 0.131
$ws.Range("E10").Value = 0.848
$ws.Range("D11").Value = 0.156
$ws.Range("E11").Value = 0.848
$ws.Range("D12").Value = 0.177
$ws.Range("E12").Value = 0.848
$ws.Range("D13").Value = 0.194
$ws.Range("E13").Value = 0.848
$ws.Range("D14").Value = 0.197
$ws.Range("E14").Value = 0.848
$ws.Range("E15").Value = 0.848
$ws.Range("D16").Value = 0.248
$ws.Range("E16").Value = 0.848
$ws.Range("D17").Value = 0.267
$ws.Range("E17").Value = 0.848
$ws.Range("A18").Formula = '=HYPERLINK("pathways/Arginine_and_Proline_Metabolism.csv","Arginine and Proline Metabolism")'
$ws.Range("B18").Value = 20
$ws.Range("D18").Value = 0.282
$ws.Range("E18").Value = 0.848
$ws.Range("A19").Formula = '=HYPERLINK("pathways/Urea_Cycle.csv","Urea Cycle")'
$ws.Range("B19").Value = 14
$ws.Range("D19").Value = 0.282
$ws.Range("E19").Value = 0.848
$ws.Range("D20").Value = 0.29
$ws.Range("E20").Value = 0.848
$ws.Range("D21").Value = 0.304
$ws.Range("E21").Value = 0.848
$ws.Range("D22").Value = 0.308
$ws.Range("E22").Value = 0.848
$ws.Range("A23").Formula = '=HYPERLINK("pathways/Glutamate_Metabolism.csv","Glutamate Metabolism")'
$ws.Range("B23").Value = 14
$ws.Range("D23").Value = 0.324
$ws.Range("E23").Value = 0.848
$ws.Range("A24").Formula = '=HYPERLINK("pathways/Fatty_Acid_Biosynthesis.csv","Fatty Acid Biosynthesis")'
$ws.Range("B24").Value = 9
$ws.Range("D24").Value = 0.326
$ws.Range("E24").Value = 0.848
$ws.Range("A25").Formula = '=HYPERLINK("pathways/Glutathione_Metabolism.csv","Glutathione Metabolism")'
$ws.Range("E25").Value = 0.848
$ws.Range("A26").Formula = '=HYPERLINK("pathways/Phospholipid_Biosynthesis.csv","Phospholipid Biosynthesis")'
$ws.Range("B26").Value = 4
$ws.Range("E26").Value = 0.848
$ws.Range("D27").Value = 0.349
$ws.Range("E27").Value = 0.848
$ws.Range("D28").Value = 0.364
$ws.Range("E28").Value = 0.848
$ws.Range("D29").Value = 0.399
$ws.Range("E29").Value = 0.848
$ws.Range("D30").Value = 0.401
$ws.Range("E30").Value = 0.848
$ws.Range("E31").Value = 0.848
$ws.Range("D32").Value = 0.407
$ws.Range("E32").Value = 0.848
$ws.Range("E33").Value = 0.848
$ws.Range("D34").Value = 0.479
$ws.Range("E34").Value = 0.848
$ws.Range("D35").Value = 0.494
$ws.Range("E35").Value = 0.848
$ws.Range("D36").Value = 0.511
$ws.Range("E36").Value = 0.848
$ws.Range("D37").Value = 0.529
$ws.Range("E37").Value = 0.848
$ws.Range("E38").Value = 0.848
$ws.Range("A39").Formula = '=HYPERLINK("pathways/Fatty_acid_Metabolism.csv","Fatty acid Metabolism")'
$ws.Range("B39").Value = 3
$ws.Range("E39").Value = 0.848
$ws.Range("A40").Formula = '=HYPERLINK("pathways/Galactose_Metabolism.csv","Galactose Metabolism")'
$ws.Range("B40").Value = 6
$ws.Range("D40").Value = 0.546
$ws.Range("E40").Value = 0.848
$ws.Range("A41").Formula = '=HYPERLINK("pathways/Retinol_Metabolism.csv","Retinol Metabolism")'
$ws.Range("D41").Value = 0.58
$ws.Range("E41").Value = 0.848
$ws.Range("A42").Formula = '=HYPERLINK("pathways/Carnitine_Synthesis.csv","Carnitine Synthesis")'
$ws.Range("D42").Value = 0.581
$ws.Range("E42").Value = 0.848
$ws.Range("D43").Value = 0.617
$ws.Range("E43").Value = 0.848
$ws.Range("E44").Value = 0.848
$ws.Range("D45").Value = 0.646
$ws.Range("E45").Value = 0.848
$ws.Range("E46").Value = 0.848
$ws.Range("E47").Value = 0.848
$ws.Range("A48").Formula = '=HYPERLINK("pathways/Tryptophan_Metabolism.csv","Tryptophan Metabolism")'
$ws.Range("B48").Value = 15
$ws.Range("E48").Value = 0.848
$ws.Range("A49").Formula = '=HYPERLINK("pathways/Vitamin_B6_Metabolism.csv","Vitamin B6 Metabolism")'
$ws.Range("B49").Value = 4
$ws.Range("D49").Value = 0.69
$ws.Range("E49").Value = 0.848
$ws.Range("A50").Formula = '=HYPERLINK("pathways/Tyrosine_Metabolism.csv","Tyrosine Metabolism")'
$ws.Range("B50").Value = 13
$ws.Range("E50").Value = 0.848
$ws.Range("A51").Formula = '=HYPERLINK("pathways/Thyroid_hormone_synthesis.csv","Thyroid hormone synthesis")'
$ws.Range("B51").Value = 5
$ws.Range("D51").Value = 0.708
$ws.Range("E51").Value = 0.848
$ws.Range("A52").Formula = '=HYPERLINK("pathways/Aspartate_Metabolism.csv","Aspartate Metabolism")'
$ws.Range("B52").Value = 14
$ws.Range("D52").Value = 0.709
$ws.Range("E52").Value = 0.848
$ws.Range("E53").Value = 0.848
$ws.Range("E54").Value = 0.848
$ws.Range("E55").Value = 0.848
$ws.Range("E56").Value = 0.848
$ws.Range("E57").Value = 0.848
$ws.Range("E58").Value = 0.848
$ws.Range("E59").Value = 0.848
$ws.Range("D60").Value = 0.78
$ws.Range("E60").Value = 0.848
$ws.Range("D61").Value = 0.784
$ws.Range("E61").Value = 0.848
$ws.Range("D62").Value = 0.799
$ws.Range("E62").Value = 0.848
$ws.Range("E63").Value = 0.848
$ws.Range("E64").Value = 0.848
$ws.Range("E65").Value = 0.848
$ws.Range("E66").Value = 0.848
$ws.Range("E67").Value = 0.848
$ws.Range("D68").Value = 0.82
$ws.Range("E68").Value = 0.848
$ws.Range("D69").Value = 0.825
$ws.Range("E69").Value = 0.848
$ws.Range("D70").Value = 0.836
$ws.Range("E70").Value = 0.848
$ws.Range("D71").Value = 0.836
$ws.Range("E71").Value = 0.848
$ws.Range("D72").Value = 0.861
$ws.Range("E72").Value = 0.861
